$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 1.85
$ws.Range("I4").Value = 5.3
$ws.Range("G6").Value = 3.05
$ws.Range("Q7").Value = 1.62
$ws.Range("F9").Value = 1.59
$ws.Range("G9").Value = 1.67
$ws.Range("H9").Value = 5.4
$ws.Range("K9").Value = 5.3
$ws.Range("P9").Value = 2.48
$ws.Range("P10").Value = 1.73
$ws.Range("Q10").Value = 2.08
